$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.402088999999998
$ws.Range("H2").Value = 28.206267
$ws.Range("I2").Value = 0.4769398566373552
$ws.Range("J2").Value = 0.4769398566373552
$ws.Range("M2").Value = 55.783591
$ws.Range("N2").Value = 167.350773
$ws.Range("O2").Value = 0.2332214199005771
$ws.Range("P2").Value = 0.2394371967339281
$ws.Range("Q2").Value = 524.482287321599
$ws.Range("R2").Value = 4720.340585894391
$ws.Range("S2").Value = 0.1112325905721417
$ws.Range("T2").Value = 0.1141971422839299

$ws.Range("G3").Value = 9.402088999999998
$ws.Range("H3").Value = 28.206267
$ws.Range("I3").Value = 0.4769398566373552
$ws.Range("J3").Value = 0.4769398566373552
$ws.Range("O3").Value = 0.5297503589663128
$ws.Range("P3").Value = 0.5438691736537713
$ws.Range("Q3").Value = 1191.334312682496
$ws.Range("R3").Value = 10722.00881414247
$ws.Range("S3").Value = 0.2526590602589807
$ws.Range("T3").Value = 0.2593928857119066

$ws.Range("G4").Value = 9.402088999999998
$ws.Range("H4").Value = 28.206267
$ws.Range("I4").Value = 0.4769398566373552
$ws.Range("J4").Value = 0.4769398566373552
$ws.Range("M4").Value = 23.03749833333333
$ws.Range("N4").Value = 69.112495
$ws.Range("O4").Value = 0.09631574403765399
$ws.Range("P4").Value = 0.09888273454277752
$ws.Range("Q4").Value = 216.6006096673516
$ws.Range("R4").Value = 1949.405487006165
$ws.Range("S4").Value = 0.04593681715323889
$ws.Range("T4").Value = 0.04716111723674196

$ws.Range("G5").Value = 9.402088999999998
$ws.Range("H5").Value = 28.206267
$ws.Range("I5").Value = 0.4769398566373552
$ws.Range("J5").Value = 0.4769398566373552
$ws.Range("M5").Value = 18.627865
$ws.Range("N5").Value = 37.25573
$ws.Range("O5").Value = 0.07787983970082285
$ws.Range("P5").Value = 0.05330365312071852
$ws.Range("Q5").Value = 175.140844609985
$ws.Range("R5").Value = 1050.84506765991
$ws.Range("S5").Value = 0.03714399958185065
$ws.Range("T5").Value = 0.0254226366776428

$ws.Range("G6").Value = 9.402088999999998
$ws.Range("H6").Value = 28.206267
$ws.Range("I6").Value = 0.4769398566373552
$ws.Range("J6").Value = 0.4769398566373552
$ws.Range("M6").Value = 15.028766
$ws.Range("N6").Value = 45.086298
$ws.Range("O6").Value = 0.06283263739463307
$ws.Range("P6").Value = 0.06450724194880479
$ws.Range("Q6").Value = 141.301795492174
$ws.Range("R6").Value = 1271.716159429566
$ws.Range("S6").Value = 0.02996738907114322
$ws.Range("T6").Value = 0.03076607472713415

$ws.Range("G7").Value = 10.311275
$ws.Range("H7").Value = 30.933825
$ws.Range("I7").Value = 0.5230601433626448
$ws.Range("J7").Value = 0.5230601433626447
$ws.Range("M7").Value = 55.783591
$ws.Range("N7").Value = 167.350773
$ws.Range("O7").Value = 0.2332214199005771
$ws.Range("P7").Value = 0.2394371967339281
$ws.Range("Q7").Value = 575.199947288525
$ws.Range("R7").Value = 5176.799525596725
$ws.Range("S7").Value = 0.1219888293284355
$ws.Range("T7").Value = 0.1252400544499982

$ws.Range("G8").Value = 10.311275
$ws.Range("H8").Value = 30.933825
$ws.Range("I8").Value = 0.5230601433626448
$ws.Range("J8").Value = 0.5230601433626447
$ws.Range("O8").Value = 0.5297503589663128
$ws.Range("P8").Value = 0.5438691736537713
$ws.Range("Q8").Value = 1306.536846758758
$ws.Range("R8").Value = 11758.83162082883
$ws.Range("S8").Value = 0.2770912987073322
$ws.Range("T8").Value = 0.2844762879418647

$ws.Range("G9").Value = 10.311275
$ws.Range("H9").Value = 30.933825
$ws.Range("I9").Value = 0.5230601433626448
$ws.Range("J9").Value = 0.5230601433626447
$ws.Range("M9").Value = 23.03749833333333
$ws.Range("N9").Value = 69.112495
$ws.Range("O9").Value = 0.09631574403765399
$ws.Range("P9").Value = 0.09888273454277752
$ws.Range("Q9").Value = 237.5459806270416
$ws.Range("R9").Value = 2137.913825643375
$ws.Range("S9").Value = 0.0503789268844151
$ws.Range("T9").Value = 0.05172161730603555

$ws.Range("G10").Value = 10.311275
$ws.Range("H10").Value = 30.933825
$ws.Range("I10").Value = 0.5230601433626448
$ws.Range("J10").Value = 0.5230601433626447
$ws.Range("M10").Value = 18.627865
$ws.Range("N10").Value = 37.25573
$ws.Range("O10").Value = 0.07787983970082285
$ws.Range("P10").Value = 0.05330365312071852
$ws.Range("Q10").Value = 192.077038677875
$ws.Range("R10").Value = 1152.46223206725
$ws.Range("S10").Value = 0.0407358401189722
$ws.Range("T10").Value = 0.02788101644307571

$ws.Range("G11").Value = 10.311275
$ws.Range("H11").Value = 30.933825
$ws.Range("I11").Value = 0.5230601433626448
$ws.Range("J11").Value = 0.5230601433626447
$ws.Range("M11").Value = 15.028766
$ws.Range("N11").Value = 45.086298
$ws.Range("O11").Value = 0.06283263739463307
$ws.Range("P11").Value = 0.06450724194880479
$ws.Range("Q11").Value = 154.96573913665
$ws.Range("R11").Value = 1394.69165222985
$ws.Range("S11").Value = 0.03286524832348985
$ws.Range("T11").Value = 0.03374116722167064
